$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New customer rows (4 and 5) - category-import feature adds two more sample rows.
$ws.Range("C4").Value = "Sus"
$ws.Range("D4").Value = "Tefy"
$ws.Range("C5").Value = "F GB h"
$ws.Range("D5").Value = "Gdgh"

# New eMail hyperlinks for the two new rows (same pattern as the existing F2/F3 mailto links).
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:Ghdc@gcfdb.fh", "", "", "Ghdc@gcfdb.fh")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:Ggsd@jhvt.yfg", "", "", "Ggsd@jhvt.yfg")

# New "Kategorie" column values for all four data rows.
$ws.Range("H2").Value = "Apple|Microsoft"
$ws.Range("H4").Value = "Dell|Apple"
$ws.Range("H3").Value = "Microsoft"
$ws.Range("H5").Value = "Dell|Microsoft Azure"

# Move the active selection like the author left it after editing.
$ws.Range("H7").Select()
